$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 5.210353775561958
$ws.Range("D2").Value = 3.450527141762511
$ws.Range("E2").Value = 10.93110129502653
$ws.Range("F2").Value = 59.8722537255562
$ws.Range("G2").Value = 3.822713803814603
$ws.Range("I2").Value = 45.682131256778
$ws.Range("J2").Value = 10.83401815670579
$ws.Range("K2").Value = 24.42152460772722
$ws.Range("L2").Value = 9.982457703458961
$ws.Range("M2").Value = 23.15127998607431
$ws.Range("C3").Value = 5.211240742280411
$ws.Range("D3").Value = 3.444814729235012
$ws.Range("E3").Value = 10.94966520184082
$ws.Range("F3").Value = 59.87336905896592
$ws.Range("G3").Value = 3.826767815703504
$ws.Range("I3").Value = 45.70349656751561
$ws.Range("J3").Value = 10.85636376982253
$ws.Range("K3").Value = 24.26149337509169
$ws.Range("L3").Value = 10.00018220333386
$ws.Range("M3").Value = 23.11722727482594
$ws.Range("C4").Value = 5.212374089271477
$ws.Range("D4").Value = 3.441300861817885
$ws.Range("E4").Value = 10.96201001890031
$ws.Range("F4").Value = 59.88653268077849
$ws.Range("G4").Value = 3.829385371300183
$ws.Range("I4").Value = 45.72570147766792
$ws.Range("J4").Value = 10.87094858669858
$ws.Range("K4").Value = 24.16886140522782
$ws.Range("L4").Value = 10.01172467129171
$ws.Range("M4").Value = 23.10091315669994
$ws.Range("C5").Value = 5.212984902565863
$ws.Range("D5").Value = 3.439867211489632
$ws.Range("E5").Value = 10.96727920275144
$ws.Range("F5").Value = 59.89502715709324
$ws.Range("G5").Value = 3.830484452589298
$ws.Range("I5").Value = 45.73702863975723
$ws.Range("J5").Value = 10.8771098840965
$ws.Range("K5").Value = 24.13256222082591
$ws.Range("L5").Value = 10.01659464772559
$ws.Range("M5").Value = 23.09542574662721
$ws.Range("C6").Value = 5.213095351323957
$ws.Range("D6").Value = 3.439629048049298
$ws.Range("E6").Value = 10.9681685723997
$ws.Range("F6").Value = 59.89662645855683
$ws.Range("G6").Value = 3.83066891511907
$ws.Range("I6").Value = 45.73904691621627
$ws.Range("J6").Value = 10.87814613542522
$ws.Range("K6").Value = 24.12662323452361
$ws.Range("L6").Value = 10.01741336338042
$ws.Range("M6").Value = 23.09458481424198
$ws.Range("C7").Value = 5.212381722541632
$ws.Range("D7").Value = 3.441281534189074
$ws.Range("E7").Value = 10.96208011423888
$ws.Range("F7").Value = 59.88663457783782
$ws.Range("G7").Value = 3.829400062520004
$ws.Range("I7").Value = 45.72584502358313
$ws.Range("J7").Value = 10.87103079735424
$ws.Range("K7").Value = 24.16836595146412
$ws.Range("L7").Value = 10.01178967541492
$ws.Range("M7").Value = 23.10083444532227
$ws.Range("C8").Value = 5.210537964274008
$ws.Range("D8").Value = 3.4485587681115
$ws.Range("E8").Value = 10.93730602219883
$ws.Range("F8").Value = 59.87004412979686
$ws.Range("G8").Value = 3.824085058083996
$ws.Range("I8").Value = 45.68760851492746
$ws.Range("J8").Value = 10.8415437790844
$ws.Range("K8").Value = 24.36519885334603
$ws.Range("L8").Value = 9.988432507792352
$ws.Range("M8").Value = 23.13858853676961
$ws.Range("C9").Value = 5.211553090639605
$ws.Range("D9").Value = 3.462784137234998
$ws.Range("E9").Value = 10.89620929394223
$ws.Range("F9").Value = 59.93680835910526
$ws.Range("G9").Value = 3.814675253905113
$ws.Range("I9").Value = 45.68497704840991
$ws.Range("J9").Value = 10.79055792857102
$ws.Range("K9").Value = 24.79424185984761
$ws.Range("L9").Value = 9.947840890349225
$ws.Range("M9").Value = 23.24880780221938
$ws.Range("C10").Value = 5.215066859548481
$ws.Range("D10").Value = 3.473217158810878
$ws.Range("E10").Value = 10.87054512237487
$ws.Range("F10").Value = 60.04669562556217
$ws.Range("G10").Value = 3.808371381346531
$ws.Range("I10").Value = 45.72744672495444
$ws.Range("J10").Value = 10.75723757667303
$ws.Range("K10").Value = 25.13342158672401
$ws.Range("L10").Value = 9.921165592951098
$ws.Range("M10").Value = 23.35143505343438
$ws.Range("C11").Value = 5.217254298672336
$ws.Range("D11").Value = 3.477961594618162
$ws.Range("E11").Value = 10.85984629960447
$ws.Range("F11").Value = 60.10992614449798
$ws.Range("G11").Value = 3.805634229278269
$ws.Range("I11").Value = 45.75645245511228
$ws.Range("J11").Value = 10.74297175209024
$ws.Range("K11").Value = 25.29240337700695
$ws.Range("L11").Value = 9.909707405561059
$ws.Range("M11").Value = 23.40272028833078
$ws.Range("C12").Value = 5.218166260582795
$ws.Range("D12").Value = 3.479758123105212
$ws.Range("E12").Value = 10.85593470739314
$ws.Range("F12").Value = 60.13577331821378
$ws.Range("G12").Value = 3.804616377032661
$ws.Range("I12").Value = 45.7688306308517
$ws.Range("J12").Value = 10.73769741919509
$ws.Range("K12").Value = 25.35323364534144
$ws.Range("L12").Value = 9.905465288832843
$ws.Range("M12").Value = 23.42279202282653
$ws.Range("C13").Value = 5.217966151271463
$ws.Range("D13").Value = 3.479371210044909
$ws.Range("E13").Value = 10.85677092848309
$ws.Range("F13").Value = 60.13012205984209
$ws.Range("G13").Value = 3.804834762255851
$ws.Range("I13").Value = 45.76610274337351
$ws.Range("J13").Value = 10.73882766256687
$ws.Range("K13").Value = 25.34010561273706
$ws.Range("L13").Value = 9.906374604029414
$ws.Range("M13").Value = 23.4184404282249
$ws.Range("C14").Value = 5.217327656549588
$ws.Range("D14").Value = 3.478109396982606
$ws.Range("E14").Value = 10.85952169132679
$ws.Range("F14").Value = 60.11201447760391
$ws.Range("G14").Value = 3.805550116914803
$ws.Range("I14").Value = 45.75744286579123
$ws.Range("J14").Value = 10.74253527007197
$ws.Range("K14").Value = 25.2973956185111
$ws.Range("L14").Value = 9.909356465524695
$ws.Range("M14").Value = 23.40435861398228
$ws.Range("C15").Value = 5.216947418515915
$ws.Range("D15").Value = 3.47733649317285
$ws.Range("E15").Value = 10.86122480629435
$ws.Range("F15").Value = 60.10117084443421
$ws.Range("G15").Value = 3.805990717256898
$ws.Range("I15").Value = 45.75232005033322
$ws.Range("J15").Value = 10.74482292125712
$ws.Range("K15").Value = 25.27131480547907
$ws.Range("L15").Value = 9.911195541677785
$ws.Range("M15").Value = 23.39581757434491
$ws.Range("C16").Value = 5.214935666402508
$ws.Range("D16").Value = 3.472907078723485
$ws.Range("E16").Value = 10.87126390956282
$ws.Range("F16").Value = 60.04282977960229
$ws.Range("G16").Value = 3.808552876663381
$ws.Range("I16").Value = 45.72574619243635
$ws.Range("J16").Value = 10.75818779137573
$ws.Range("K16").Value = 25.12312210340685
$ws.Range("L16").Value = 9.92192798884825
$ws.Range("M16").Value = 23.34817513606562
$ws.Range("C17").Value = 5.213851624274417
$ws.Range("D17").Value = 3.470189467059274
$ws.Range("E17").Value = 10.87767216542769
$ws.Range("F17").Value = 60.01043088867449
$ws.Range("G17").Value = 3.81015802034438
$ws.Range("I17").Value = 45.7119262233429
$ws.Range("J17").Value = 10.76661481658731
$ws.Range("K17").Value = 25.03337830524083
$ws.Range("L17").Value = 9.928684964605099
$ws.Range("M17").Value = 23.3201190290288
$ws.Range("C18").Value = 5.213283617618595
$ws.Range("D18").Value = 3.468626209292959
$ws.Range("E18").Value = 10.88144991548701
$ws.Range("F18").Value = 59.99304250476685
$ws.Range("G18").Value = 3.811093547986154
$ws.Range("I18").Value = 45.70488901919416
$ws.Range("J18").Value = 10.77154578050282
$ws.Range("K18").Value = 24.98220410140492
$ws.Range("L18").Value = 9.932635103242573
$ws.Range("M18").Value = 23.30441549724481
$ws.Range("C19").Value = 5.213100865821709
$ws.Range("D19").Value = 3.468096885443205
$ws.Range("E19").Value = 10.88274479491272
$ws.Range("F19").Value = 59.98736920489053
$ws.Range("G19").Value = 3.811412416286882
$ws.Range("I19").Value = 45.70266286105503
$ws.Range("J19").Value = 10.77322975326986
$ws.Range("K19").Value = 24.96495506777554
$ws.Range("L19").Value = 9.93398350700361
$ws.Range("M19").Value = 23.29917332320986
$ws.Range("C20").Value = 5.213961285484084
$ws.Range("D20").Value = 3.470478776940964
$ws.Range("E20").Value = 10.87698048907228
$ws.Range("F20").Value = 60.01375078955189
$ws.Range("G20").Value = 3.809985878808201
$ws.Range("I20").Value = 45.71330301342027
$ws.Range("J20").Value = 10.76570905862071
$ws.Range("K20").Value = 25.04288606277357
$ws.Range("L20").Value = 9.927959082952222
$ws.Range("M20").Value = 23.32306084580533
$ws.Range("C21").Value = 5.217512936985655
$ws.Range("D21").Value = 3.478480023235313
$ws.Range("E21").Value = 10.85870993505657
$ws.Range("F21").Value = 60.11728147991566
$ws.Range("G21").Value = 3.805339494764336
$ws.Range("I21").Value = 45.75994863329897
$ws.Range("J21").Value = 10.74144278997641
$ws.Range("K21").Value = 25.30992393109356
$ws.Range("L21").Value = 9.908477995488726
$ws.Range("M21").Value = 23.40847719838621
$ws.Range("C22").Value = 5.220321072042624
$ws.Range("D22").Value = 3.483708868077945
$ws.Range("E22").Value = 10.84758380917927
$ws.Range("F22").Value = 60.19603601902589
$ws.Range("G22").Value = 3.802411452888244
$ws.Range("I22").Value = 45.79856143118973
$ws.Range("J22").Value = 10.72632825737393
$ws.Range("K22").Value = 25.48808278050327
$ws.Range("L22").Value = 9.896310295942174
$ws.Range("M22").Value = 23.46809170954537
$ws.Range("C23").Value = 5.218778123398041
$ws.Range("D23").Value = 3.480918125592152
$ws.Range("E23").Value = 10.85344765431399
$ws.Range("F23").Value = 60.15298920445264
$ws.Range("G23").Value = 3.803964302873708
$ws.Range("I23").Value = 45.77720919943354
$ws.Range("J23").Value = 10.73432714580449
$ws.Range("K23").Value = 25.39267919537489
$ws.Range("L23").Value = 9.902752934052442
$ws.Range("M23").Value = 23.43593112824383
$ws.Range("C24").Value = 5.213911535560995
$ws.Range("D24").Value = 3.470347982653645
$ws.Range("E24").Value = 10.87729290449993
$ws.Range("F24").Value = 60.01224600568115
$ws.Range("G24").Value = 3.810063664401439
$ws.Range("I24").Value = 45.71267773828437
$ws.Range("J24").Value = 10.76611828341321
$ws.Range("K24").Value = 25.03858629167325
$ws.Range("L24").Value = 9.928287050111148
$ws.Range("M24").Value = 23.32172952184706
$ws.Range("C25").Value = 5.2107883215903
$ws.Range("D25").Value = 3.458939713302871
$ws.Range("E25").Value = 10.90652918561636
$ws.Range("F25").Value = 59.90807855162333
$ws.Range("G25").Value = 3.817113243794968
$ws.Range("I25").Value = 45.6779088082789
$ws.Range("J25").Value = 10.80362197886674
$ws.Range("K25").Value = 24.67379361160249
$ws.Range("L25").Value = 9.958267145464603
$ws.Range("M25").Value = 23.21515592600442
